# Updates the "loading_percent" results sheet (case with 380 kV) with the
# recomputed loading percentages for rows 2-25, columns C-H and K-O.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 13.18025290848255
$ws.Range("D2").Value = 5.092998200233177
$ws.Range("E2").Value = 13.61565462538543
$ws.Range("F2").Value = 25.14215664457747
$ws.Range("G2").Value = 30.28802097968242
$ws.Range("H2").Value = 14.64776593385198
$ws.Range("K2").Value = 12.598545853874
$ws.Range("L2").Value = 9.309667080691684
$ws.Range("M2").Value = 16.30038046424097
$ws.Range("N2").Value = 17.6442715975552
$ws.Range("O2").Value = 22.50764160354123

$ws.Range("C3").Value = 13.14587934746885
$ws.Range("D3").Value = 5.045562305372663
$ws.Range("E3").Value = 13.63393945377854
$ws.Range("F3").Value = 25.15756124184188
$ws.Range("G3").Value = 30.30475339533857
$ws.Range("H3").Value = 14.69267664512524
$ws.Range("K3").Value = 12.18273426213017
$ws.Range("L3").Value = 9.332724359298773
$ws.Range("M3").Value = 16.13825030704243
$ws.Range("N3").Value = 17.68154504779005
$ws.Range("O3").Value = 22.56972642742681

$ws.Range("C4").Value = 13.12779032418337
$ws.Range("D4").Value = 5.015788849116166
$ws.Range("E4").Value = 13.64772829990681
$ws.Range("F4").Value = 25.17409836901198
$ws.Range("G4").Value = 30.32540677633289
$ws.Range("H4").Value = 14.7228128199231
$ws.Range("K4").Value = 11.92106831689177
$ws.Range("L4").Value = 9.34789895055912
$ws.Range("M4").Value = 16.04021432440845
$ws.Range("N4").Value = 17.70616152429199
$ws.Range("O4").Value = 22.61311715221149

$ws.Range("C5").Value = 13.12118320936953
$ws.Range("D5").Value = 5.00349836830542
$ws.Range("E5").Value = 13.65399126915129
$ws.Range("F5").Value = 25.18261504076073
$ws.Range("G5").Value = 30.33642668057094
$ws.Range("H5").Value = 14.73573694176677
$ws.Range("K5").Value = 11.81299467292969
$ws.Range("L5").Value = 9.354338860266537
$ws.Range("M5").Value = 16.00068229991091
$ws.Range("N5").Value = 17.71662879650859
$ws.Range("O5").Value = 22.63212108709371

$ws.Range("C6").Value = 13.12013241907518
$ws.Range("D6").Value = 5.001448183474101
$ws.Range("E6").Value = 13.6550701085268
$ws.Range("F6").Value = 25.18413650508147
$ws.Range("G6").Value = 30.33841354078103
$ws.Range("H6").Value = 14.73792182478785
$ws.Range("K6").Value = 11.79496680757942
$ws.Range("L6").Value = 9.35542368222403
$ws.Range("M6").Value = 15.99414440813422
$ws.Range("N6").Value = 17.71839322544285
$ws.Range("O6").Value = 22.6353564089419

$ws.Range("C7").Value = 13.12769811659597
$ws.Range("D7").Value = 5.015623726292868
$ws.Range("E7").Value = 13.64781015788162
$ws.Range("F7").Value = 25.17420603403008
$ws.Range("G7").Value = 30.3255448633717
$ws.Range("H7").Value = 14.72298451485725
$ws.Range("K7").Value = 11.91961642359524
$ws.Range("L7").Value = 9.34798476389699
$ws.Range("M7").Value = 16.03967943772003
$ws.Range("N7").Value = 17.70630092384328
$ws.Range("O7").Value = 22.61336809845274

$ws.Range("C8").Value = 13.1677785423668
$ws.Range("D8").Value = 5.076779867739202
$ws.Range("E8").Value = 13.62142737561364
$ws.Range("F8").Value = 25.14599789363706
$ws.Range("G8").Value = 30.29163268656142
$ws.Range("H8").Value = 14.66271932406791
$ws.Range("K8").Value = 12.45658235595809
$ws.Range("L8").Value = 9.31740627915994
$ws.Range("M8").Value = 16.24418710121724
$ws.Range("N8").Value = 17.65676476802142
$ws.Range("O8").Value = 22.52795258864746

$ws.Range("C9").Value = 13.27002125398887
$ws.Range("D9").Value = 5.191341332432606
$ws.Range("E9").Value = 13.59002665573663
$ws.Range("F9").Value = 25.14690598591412
$ws.Range("G9").Value = 30.30768990758667
$ws.Range("H9").Value = 14.56488257135022
$ws.Range("K9").Value = 13.45298993406675
$ws.Range("L9").Value = 9.26549945243876
$ws.Range("M9").Value = 16.65552670538715
$ws.Range("N9").Value = 17.57332274561353
$ws.Range("O9").Value = 22.40241573807503

$ws.Range("C10").Value = 13.3591164832932
$ws.Range("D10").Value = 5.271943402274236
$ws.Range("E10").Value = 13.57935747774388
$ws.Range("F10").Value = 25.18186142998796
$ws.Range("G10").Value = 30.36995703677539
$ws.Range("H10").Value = 14.50543586205771
$ws.Range("K10").Value = 14.14326143824805
$ws.Range("L10").Value = 9.232255569813729
$ws.Range("M10").Value = 16.96171464650694
$ws.Range("N10").Value = 17.52032582338757
$ws.Range("O10").Value = 22.33595014045624

$ws.Range("C11").Value = 13.40257465817503
$ws.Range("D11").Value = 5.307776511922943
$ws.Range("E11").Value = 13.57719380513918
$ws.Range("F11").Value = 25.20518763457832
$ws.Range("G11").Value = 30.4092271189037
$ws.Range("H11").Value = 14.48109846157129
$ws.Range("K11").Value = 14.44692296340724
$ws.Range("L11").Value = 9.21819015038613
$ws.Range("M11").Value = 17.10137981561798
$ws.Range("N11").Value = 17.49801130019983
$ws.Range("O11").Value = 22.31134247107963

$ws.Range("C12").Value = 13.41944194350207
$ws.Range("D12").Value = 5.321220792674589
$ws.Range("E12").Value = 13.57676068479
$ws.Range("F12").Value = 25.21508449981924
$ws.Range("G12").Value = 30.42566638741548
$ws.Range("H12").Value = 14.47227204822001
$ws.Range("K12").Value = 14.56033247572602
$ws.Range("L12").Value = 9.213015684449408
$ws.Range("M12").Value = 17.15428094232846
$ws.Range("N12").Value = 17.48981868491729
$ws.Range("O12").Value = 22.3028358095231

$ws.Range("C13").Value = 13.41579118314999
$ws.Range("D13").Value = 5.318330963231396
$ws.Range("E13").Value = 13.57683679917016
$ws.Range("F13").Value = 25.21290579845746
$ws.Range("G13").Value = 30.42205624364478
$ws.Range("H13").Value = 14.47415563560969
$ws.Range("K13").Value = 14.53597936300369
$ws.Range("L13").Value = 9.214123350700412
$ws.Range("M13").Value = 17.14288781696249
$ws.Range("N13").Value = 17.4915716717619
$ws.Range("O13").Value = 22.30463173371163

$ws.Range("C14").Value = 13.4039541728004
$ws.Range("D14").Value = 5.30888511387583
$ws.Range("E14").Value = 13.57715043570833
$ws.Range("F14").Value = 25.20598056363804
$ws.Range("G14").Value = 30.41054818979397
$ws.Range("H14").Value = 14.4803644949438
$ws.Range("K14").Value = 14.45628540058259
$ws.Range("L14").Value = 9.217761403006881
$ws.Range("M14").Value = 17.10573199928455
$ws.Range("N14").Value = 17.49733213314746
$ws.Range("O14").Value = 22.31062634172465

$ws.Range("C15").Value = 13.39675681476458
$ws.Range("D15").Value = 5.303082837582628
$ws.Range("E15").Value = 13.5773928234723
$ws.Range("F15").Value = 25.2018770497246
$ws.Range("G15").Value = 30.40370324962596
$ws.Range("H15").Value = 14.48421836024375
$ws.Range("K15").Value = 14.4072621519539
$ws.Range("L15").Value = 9.22000957703764
$ws.Range("M15").Value = 17.08297345282557
$ws.Range("N15").Value = 17.5008940877792
$ws.Range("O15").Value = 22.31440399071976

$ws.Range("C16").Value = 13.3563344346246
$ws.Range("D16").Value = 5.269584416974429
$ws.Range("E16").Value = 13.57955295974855
$ws.Range("F16").Value = 25.18048610466417
$ws.Range("G16").Value = 30.36761047781919
$ws.Range("H16").Value = 14.50708085832152
$ws.Range("K16").Value = 14.12320007023337
$ws.Range("L16").Value = 9.233196034133488
$ws.Range("M16").Value = 16.95259126121426
$ws.Range("N16").Value = 17.52182018159777
$ws.Range("O16").Value = 22.33767177281983

$ws.Range("C17").Value = 13.33227938725973
$ws.Range("D17").Value = 5.248817131319722
$ws.Range("E17").Value = 13.58156670709455
$ws.Range("F17").Value = 25.16926270568514
$ws.Range("G17").Value = 30.34826885565964
$ws.Range("H17").Value = 14.52179958372709
$ws.Range("K17").Value = 13.94621772828092
$ws.Range("L17").Value = 9.241556143523187
$ws.Range("M17").Value = 16.87267333302603
$ws.Range("N17").Value = 17.53511674954669
$ws.Range("O17").Value = 22.35338907314303

$ws.Range("C18").Value = 13.3187198566426
$ws.Range("D18").Value = 5.236794213884569
$ws.Range("E18").Value = 13.58297818362175
$ws.Range("F18").Value = 25.16350647359184
$ws.Range("G18").Value = 30.33817467040003
$ws.Range("H18").Value = 14.53051999635302
$ws.Range("K18").Value = 13.84345312366484
$ws.Range("L18").Value = 9.246464198434602
$ws.Range("M18").Value = 16.82674476122674
$ws.Range("N18").Value = 17.54293348081357
$ws.Range("O18").Value = 22.36295884224394

$ws.Range("C19").Value = 13.31417659479457
$ws.Range("D19").Value = 5.232710201303061
$ws.Range("E19").Value = 13.58349959077386
$ws.Range("F19").Value = 25.16167769554829
$ws.Range("G19").Value = 30.33493408264437
$ws.Range("H19").Value = 14.53351628492771
$ws.Range("K19").Value = 13.80849536612057
$ws.Range("L19").Value = 9.248143085401187
$ws.Range("M19").Value = 16.81120195342621
$ws.Range("N19").Value = 17.54560911755265
$ws.Range("O19").Value = 22.36628988028818

$ws.Range("C20").Value = 13.33481156621709
$ws.Range("D20").Value = 5.251035967460727
$ws.Range("E20").Value = 13.58132613622844
$ws.Range("F20").Value = 25.1703851188318
$ws.Range("G20").Value = 30.35022117146954
$ws.Range("H20").Value = 14.52020639641657
$ws.Range("K20").Value = 13.96515879929528
$ws.Range("L20").Value = 9.240655895284473
$ws.Range("M20").Value = 16.88117709124508
$ws.Range("N20").Value = 17.53368383114611
$ws.Range("O20").Value = 22.35166111001377

$ws.Range("C21").Value = 13.40741993231738
$ws.Range("D21").Value = 5.311663019846302
$ws.Range("E21").Value = 13.57704783680984
$ws.Range("F21").Value = 25.20798584086368
$ws.Range("G21").Value = 30.41388586769652
$ws.Range("H21").Value = 14.47853022310461
$ws.Range("K21").Value = 14.47973702893467
$ws.Range("L21").Value = 9.21668870092291
$ws.Range("M21").Value = 17.11664554004732
$ws.Range("N21").Value = 17.49563316444565
$ws.Range("O21").Value = 22.30884353195885

$ws.Range("C22").Value = 13.45726178087225
$ws.Range("D22").Value = 5.35055577539341
$ws.Range("E22").Value = 13.57650247639861
$ws.Range("F22").Value = 25.23875767857471
$ws.Range("G22").Value = 30.4646330277489
$ws.Range("H22").Value = 14.45356377749241
$ws.Range("K22").Value = 14.80679567340735
$ws.Range("L22").Value = 9.201909466194168
$ws.Range("M22").Value = 17.27059447578413
$ws.Range("N22").Value = 17.47226498273545
$ws.Range("O22").Value = 22.2855920631859

$ws.Range("C23").Value = 13.4304453645631
$ws.Range("D23").Value = 5.329866482077368
$ws.Range("E23").Value = 13.57658784102106
$ws.Range("F23").Value = 25.22176870526713
$ws.Range("G23").Value = 30.43671445967318
$ws.Range("H23").Value = 14.46668081213539
$ws.Range("K23").Value = 14.63311236969389
$ws.Range("L23").Value = 9.209716551740346
$ws.Range("M23").Value = 17.18843740961072
$ws.Range("N23").Value = 17.48459994456461
$ws.Range("O23").Value = 22.29756807526728

$ws.Range("C24").Value = 13.3336659262007
$ws.Range("D24").Value = 5.250033091419265
$ws.Range("E24").Value = 13.58143410788518
$ws.Range("F24").Value = 25.1698755067875
$ws.Range("G24").Value = 30.34933533452715
$ws.Range("H24").Value = 14.52092587146793
$ws.Range("K24").Value = 13.95659869573141
$ws.Range("L24").Value = 9.241062580653837
$ws.Range("M24").Value = 16.87733248505439
$ws.Range("N24").Value = 17.53433111662201
$ws.Range("O24").Value = 22.35244065992657

$ws.Range("C25").Value = 13.2398726362755
$ws.Range("D25").Value = 5.16095471311756
$ws.Range("E25").Value = 13.59634320674886
$ws.Range("F25").Value = 25.1406352561333
$ws.Range("G25").Value = 30.29448442990602
$ws.Range("H25").Value = 14.58916871533265
$ws.Range("K25").Value = 13.19029649648186
$ws.Range("L25").Value = 9.27868109675836
$ws.Range("M25").Value = 16.65552670538715
$ws.Range("N25").Value = 17.59443401695872
$ws.Range("O25").Value = 22.43186428754839
